$d = $word.ActiveDocument

$replacements = @(
    @("714×3=", "304×3="),
    @("575×3=", "947×4="),
    @("137×2=", "989×3="),
    @("375×4=", "774×4="),
    @("456×5=", "521×3="),
    @("700×9=", "982×2="),
    @("180×7=", "627×5="),
    @("155×6=", "801×2="),
    @("868×7=", "122×7="),
    @("854×6=", "468×2="),
    @("110×8=", "805×9="),
    @("649×8=", "712×9="),
    @("481×3=", "319×9="),
    @("447×6=", "982×2="),
    @("525×4=", "128×6="),
    @("274×5=", "163×9="),
    @("305×7=", "713×4="),
    @("187×3=", "726×5="),
    @("967×7=", "477×2="),
    @("778×4=", "756×8="),
    @("774×7=", "323×7="),
    @("422×8=", "304×5="),
    @("926×6=", "148×4="),
    @("427×8=", "237×7="),
    @("910×4=", "106×9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

$d.Save()
